# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# For numeric-looking "Price" cells a leading apostrophe is used so Excel
# keeps storing them as text (matching the rest of that column, e.g.
# "38.627.28") instead of auto-converting to a Number; the Style reset
# afterwards clears the resulting quote-prefix formatting so no stray
# number format is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.627.28"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").Value = "2.091.24"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'228.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").Value = "'0.0841"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "2.395.23"
$ws.Range("E12").Value = "  +2.67%  "
$ws.Range("D13").Value = "'14.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").Value = "'22.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.14%  "
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "'5.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.05%  "
$ws.Range("D17").Value = "2.117.07"
$ws.Range("E17").Value = "  +3.38%  "
$ws.Range("D18").Value = "38.542.95"
$ws.Range("E18").Value = "  +2.24%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'6.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'70.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("D21").Value = "0.0₃0835"
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("D22").Value = "'226.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'2.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").Value = "'170.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").Value = "'9.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  +4.51%  "
$ws.Range("D29").Value = "'19.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("E30").Value = "  +7.77%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  +4.50%  "
$ws.Range("E33").Value = "  +6.31%  "
$ws.Range("E34").Value = "  +2.38%  "
$ws.Range("D35").Value = "'0.0607"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").Value = "'6.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("D37").Value = "'2.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("D38").Value = "'3.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.09%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "'18.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").Value = "1.544.61"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").Value = "'99.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.96%  "
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").Value = "'0.0914"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("E46").Value = "  +9.58%  "
$ws.Range("D47").Value = "'4.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "2.287.08"
$ws.Range("E51").Value = "  +2.92%  "
